# Adds a new "Knärot" section (heading, body paragraphs, references) right
# after the "BILAGA 1 - Fridlysta arter" title paragraph, and bumps the date
# in the first-page header from 2023-09-13 to 2023-09-15.

$d = $word.ActiveDocument

function Insert-Run([int]$pos, [string]$text, [bool]$italic) {
    $rng = $d.Range($pos, $pos)
    $rng.InsertAfter($text)
    $endPos = $pos + $text.Length
    if ($italic) {
        $fmtRng = $d.Range($pos, $endPos)
        $fmtRng.Font.Italic = 1
    }
    return $endPos
}

# Step 1: create all the empty paragraphs first (in document order), each with
# its target style, while the "typing position" formatting is still clean.
# This avoids italic formatting leaking from the end of one paragraph into the
# start of the next paragraph that is created afterwards.
$paras = New-Object System.Collections.ArrayList
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Heading1")
[void]$paras.Add($newPara)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Normal")
[void]$paras.Add($newPara)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Normal")
[void]$paras.Add($newPara)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Normal")
[void]$paras.Add($newPara)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Normal")
[void]$paras.Add($newPara)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Normal")
[void]$paras.Add($newPara)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Heading2")
[void]$paras.Add($newPara)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Normal")
[void]$paras.Add($newPara)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Normal")
[void]$paras.Add($newPara)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Normal")
[void]$paras.Add($newPara)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Normal")
[void]$paras.Add($newPara)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Normal")
[void]$paras.Add($newPara)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles.Item("Normal")
[void]$paras.Add($newPara)

# Step 2: fill the paragraphs with their runs, starting from the LAST
# paragraph and working backwards. Because each paragraph's own runs are
# still written left-to-right, the text/formatting comes out correct, but we
# never type into a paragraph that has already-typed content after it, so a
# trailing italic run can't leak into a subsequent paragraph.

# paragraph 12
$pos = $paras[12].Range.Start
$pos = Insert-Run $pos "SLU Artdatabanken, 2021. " $false
$pos = Insert-Run $pos "Artfaktablad. Naturvård – artfakta. " $true
$pos = Insert-Run $pos "SLU Artdatabanken, Uppsala " $false

# paragraph 11
$pos = $paras[11].Range.Start
$pos = Insert-Run $pos "Skogsstyrelsen, 2022. " $false
$pos = Insert-Run $pos "Vägledning för hänsyn till knärot. " $true
$pos = Insert-Run $pos "https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/" $false

# paragraph 10
$pos = $paras[10].Range.Start
$pos = Insert-Run $pos "Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. " $false
$pos = Insert-Run $pos "Biological legacies buffer local species extinction after logging. " $true
$pos = Insert-Run $pos "Journal of Applied Ecology. 51, 53-62." $false

# paragraph 9
$pos = $paras[9].Range.Start
$pos = Insert-Run $pos "Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. " $false
$pos = Insert-Run $pos "Interactive effects of drought and edge exposure on old-growth forest understory species. " $true
$pos = Insert-Run $pos "Landscape Ecology, 37, sid 1839-1853" $false

# paragraph 8
$pos = $paras[8].Range.Start
$pos = Insert-Run $pos "Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. " $false
$pos = Insert-Run $pos "Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. " $true
$pos = Insert-Run $pos "Ecological Applications, 22, 2049-2064 " $false

# paragraph 7
$pos = $paras[7].Range.Start
$pos = Insert-Run $pos "de Graaf M & Roberts M.R., 2009. " $false
$pos = Insert-Run $pos "Short-term response of the herbaceous layer within leave patches after harvest. " $true
$pos = Insert-Run $pos "Forest Ecology and Management 257, 1014-1025" $false

# paragraph 6
$pos = $paras[6].Range.Start
$pos = Insert-Run $pos "Referenser - knärot" $false

# paragraph 5
$pos = $paras[5].Range.Start
$pos = Insert-Run $pos "Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022)." $false

# paragraph 4
$pos = $paras[4].Range.Start
$pos = Insert-Run $pos "En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022)." $false

# paragraph 3
$pos = $paras[3].Range.Start
$pos = Insert-Run $pos "Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: " $false
$pos = Insert-Run $pos "“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”" $true

# paragraph 2
$pos = $paras[2].Range.Start
$pos = Insert-Run $pos "Samuel Johnsons doktorsavhandling " $false
$pos = Insert-Run $pos "“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“" $true
$pos = Insert-Run $pos " (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: " $false
$pos = Insert-Run $pos "“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” " $true
$pos = Insert-Run $pos "Vidare " $false
$pos = Insert-Run $pos "“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”" $true

# paragraph 1
$pos = $paras[1].Range.Start
$pos = Insert-Run $pos "Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021)." $false

# paragraph 0
$pos = $paras[0].Range.Start
$pos = Insert-Run $pos "Knärot – ekologi samt krav på livsmiljön" $false

# Step 3: bump the date in the first-page header (2023-09-13 -> 2023-09-15).
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(2)
$hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null

Write-Output "Done"
